# edit.ps1 - apply resume content updates via Word COM-interop
$d = $word.ActiveDocument
$bullet = [char]0x2022

function Insert-BulletAfter($searchText, $newText) {
    $anchor = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*$searchText*") {
            $anchor = $p
            break
        }
    }
    if ($anchor -eq $null) {
        throw "Insert-BulletAfter: anchor paragraph not found for '$searchText'"
    }
    $anchor.Range.InsertParagraphAfter()
    $newPara = $anchor.Next()
    $newPara.Range.Text = $newText
}

# 1. Update years of experience in PROFESSIONAL SUMMARY
$d.Content.Find.Execute(
    "Senior Software Engineer with 21 years of experience",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Senior Software Engineer with 15+ years of experience", 2) | Out-Null

# 2. Enhance FLEEM / Twilio bullet (Progressive Change Campaign Committee)
$d.Content.Find.Execute(
    "using Twilio API for thousands of simultaneous phone calls",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys", 2) | Out-Null

# 3. Enhance Salsa Labs CRM bullet
$d.Content.Find.Execute(
    "comprehensive geospatial analysis and reporting tools for Java-based CRM system",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands of users simultaneously", 2) | Out-Null

# 4. Enhance Salsa Labs mapping/visualization bullet
$d.Content.Find.Execute(
    "Integrated mapping and visualization tools for political campaign data analysis",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Integrated mapping and visualization tools for political campaign data analysis interfacing with Government and Activism APIs", 2) | Out-Null

# 5. New bullet after "Collaborated with political strategists..." (Salsa Labs)
Insert-BulletAfter "Collaborated with political strategists to translate geospatial requirements into technical solutions" `
    "$bullet Handled billions of records with millions of columns in high-performance CRM system"

# 6. New bullet at end of Praxis Project role (before Lake Research Partners heading)
Insert-BulletAfter "Managed technology infrastructure supporting community health initiatives across multiple countries" `
    "$bullet Architected and developed 25 Drupal sites to integrate with membership databases, activism CRMs and government agencies, under guidelines from Kellogg Foundation and Robert Wood Johnson Foundation"

# 7. New bullet at end of Lake Research Partners role
Insert-BulletAfter "Developed innovative approaches to visualizing demographic and market data for enhanced client understanding" `
    "$bullet Trained staff on building Python tooling for report generation and analysis"

# 8. Replace the EDUCATION section (heading + two degree lines) with a single bullet line
$ranges = @()
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*EDUCATION*" -or $t -like "*Master of Arts in Political Science*" -or $t -like "*Bachelor of Arts in Political Science*") {
        $ranges += ,@($p.Range.Start, $p.Range.End)
    }
}
for ($i = $ranges.Count - 1; $i -ge 0; $i--) {
    $s = $ranges[$i][0]
    $e = $ranges[$i][1]
    $d.Range($s, $e).Delete()
}

Insert-BulletAfter "Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL" `
    "$bullet Trained staff on PHP/MySQL for data analysis and reporting systems"

Write-Output "edit.ps1 completed"
